$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 10.56472949907906
$ws.Range("C2").Value = 6.212147121844095
$ws.Range("D2").Value = 4.559132871428442
$ws.Range("E2").Value = 3.99283691009699
$ws.Range("B3").Value = 8.374954051276625
$ws.Range("C3").Value = 4.809282045285278
$ws.Range("D3").Value = 3.651525841603914
$ws.Range("E3").Value = 3.23283331388647
$ws.Range("B4").Value = 9.361953554425824
$ws.Range("C4").Value = 5.544570533031403
$ws.Range("D4").Value = 4.19838724095648
$ws.Range("E4").Value = 3.743511689471815
$ws.Range("B5").Value = 4.940140009134441
$ws.Range("C5").Value = 2.934576195968682
$ws.Range("D5").Value = 2.42626073629948
$ws.Range("E5").Value = 2.286912720564253
$ws.Range("B6").Value = 1.426955281458697
$ws.Range("C6").Value = 0.8753565019908757
$ws.Range("D6").Value = 0.6342104514815061
$ws.Range("E6").Value = 0.5354154108338075
$ws.Range("B7").Value = 1.909465123111566
$ws.Range("C7").Value = 1.219932049956222
$ws.Range("D7").Value = 1.000518488509308
$ws.Range("E7").Value = 0.9101152687398403
$ws.Range("B8").Value = 3.162607168728365
$ws.Range("C8").Value = 1.962454242750496
$ws.Range("D8").Value = 1.398601356383271
$ws.Range("E8").Value = 1.149297017754302
$ws.Range("B9").Value = 2.414254612245162
$ws.Range("C9").Value = 1.489277673447328
$ws.Range("D9").Value = 1.032885358940224
$ws.Range("E9").Value = 0.8337038804962449
$ws.Range("B10").Value = 4.478426308014347
$ws.Range("C10").Value = 3.052325732473214
$ws.Range("D10").Value = 2.919303411800664
$ws.Range("E10").Value = 2.973400820959878
$ws.Range("B11").Value = 3.624572709325329
$ws.Range("C11").Value = 2.462184757981162
$ws.Range("D11").Value = 2.465559063346524
$ws.Range("E11").Value = 2.583480844780979
$ws.Range("B12").Value = 0.7992521905295669
$ws.Range("C12").Value = 0.6115296961620851
$ws.Range("D12").Value = 0.9315824148048609
$ws.Range("E12").Value = 1.194841888868514
$ws.Range("B13").Value = 3.085090764804561
$ws.Range("C13").Value = 2.167242324121294
$ws.Range("D13").Value = 2.301587582174857
$ws.Range("E13").Value = 2.47231477791947
